$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings (e.g. "219.09") must be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# and they lose their original text representation (trailing zeros, etc).
$textForcedCells = @('D5', 'D8', 'D9', 'D11', 'D16', 'D18', 'D19', 'D25', 'D26', 'D27', 'D31', 'D34', 'D36', 'D37', 'D40', 'D41', 'D46', 'D48', 'D50')
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.053.90'
$ws.Range('E2').Value = '  +3.41%  '

$ws.Range('D3').Value = '1.723.09'
$ws.Range('E3').Value = '  +2.56%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '219.09'
$ws.Range('E5').Value = '  +1.84%  '

$ws.Range('E6').Value = '  +0.58%  '

$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('D8').Value = '24.36'
$ws.Range('E8').Value = '  +14.25%  '

$ws.Range('D9').Value = '0.264'
$ws.Range('E9').Value = '  +3.07%  '

$ws.Range('E10').Value = '  +1.56%  '

$ws.Range('D11').Value = '0.0898'
$ws.Range('E11').Value = '  +1.58%  '

$ws.Range('D12').Value = '1.966.00'
$ws.Range('E12').Value = '  +2.55%  '

$ws.Range('D13').Value = '1.731.28'
$ws.Range('E13').Value = '  +2.53%  '

$ws.Range('E14').Value = '  +3.01%  '

$ws.Range('E15').Value = '  +4.76%  '

$ws.Range('D16').Value = '67.55'

$ws.Range('D17').Value = '28.013.03'
$ws.Range('E17').Value = '  +3.31%  '

$ws.Range('D18').Value = '242.46'
$ws.Range('E18').Value = '  +1.80%  '

$ws.Range('D19').Value = '8.01'
$ws.Range('E19').Value = '  -1.57%  '

$ws.Range('E20').Value = '  +1.10%  '

$ws.Range('E21').Value = '  -0.12%  '

$ws.Range('E22').Value = '  +2.43%  '

$ws.Range('E23').Value = '  +1.99%  '

$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('D25').Value = '148.75'
$ws.Range('E25').Value = '  +1.36%  '

$ws.Range('D26').Value = '7.48'
$ws.Range('E26').Value = '  +3.38%  '

$ws.Range('D27').Value = '16.72'
$ws.Range('E27').Value = '  +2.51%  '

$ws.Range('E28').Value = '  +0.93%  '

$ws.Range('E29').Value = '  +0.12%  '

$ws.Range('E30').Value = '  +1.93%  '

$ws.Range('D31').Value = '1.20'
$ws.Range('E31').Value = '  +2.00%  '

$ws.Range('E32').Value = '  +2.05%  '

$ws.Range('D33').Value = '1.492.98'
$ws.Range('E33').Value = '  -4.51%  '

$ws.Range('D34').Value = '3.27'
$ws.Range('E34').Value = '  +2.06%  '

$ws.Range('E35').Value = '  -2.36%  '

$ws.Range('D36').Value = '0.954'
$ws.Range('E36').Value = '  +1.92%  '

$ws.Range('D37').Value = '0.607'
$ws.Range('E37').Value = '  +0.62%  '

$ws.Range('E38').Value = '  +0.79%  '

$ws.Range('E39').Value = '  -0.05%  '

$ws.Range('D40').Value = '1.06'
$ws.Range('E40').Value = '  +0.54%  '

$ws.Range('D41').Value = '70.61'
$ws.Range('E41').Value = '  +2.28%  '

$ws.Range('E42').Value = '  +3.18%  '

$ws.Range('E43').Value = '  -0.10%  '

$ws.Range('D45').Value = '1.870.05'
$ws.Range('E45').Value = '  +2.37%  '

$ws.Range('D46').Value = '0.801'
$ws.Range('E46').Value = '  +2.32%  '

$ws.Range('E47').Value = '  +11.58%  '

$ws.Range('D48').Value = '91.05'
$ws.Range('E48').Value = '  +0.29%  '

$ws.Range('E49').Value = '  +4.42%  '

$ws.Range('D50').Value = '8.27'
$ws.Range('E50').Value = '  +2.34%  '

$ws.Range('E51').Value = '  +0.60%  '
